$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shared string update: "E7420" -> "E7420L" for the sample-number column (G2:G33).
#    All of those cells share the same string, so re-writing the range text
#    updates the single shared-string table entry referenced by every cell.
$ws.Range("G2:G33").Value = "E7420L"

# 2. H2:H33 - add an explicit =FALSE() formula on top of the existing boolean
#    value (major accuracy check update). Each cell gets its own formula
#    (not a shared/array formula) so every row is an independent =FALSE().
for ($r = 2; $r -le 33; $r++) {
    $ws.Range("H$r").Formula = "=FALSE()"
}

# 3. Update the view: scroll so row 3 is the top-left visible row, and move
#    the active selection from the H column to the G column (G2:G33).
$null = $ws.Range("G2:G33").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
